$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 becomes the last row of its group: reuse the existing
# "bottom border" formatting (styles 6/7, already present in styles.xml,
# e.g. row 4) so no new style entries are created.
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null

# Row 18 starts a new group: reuse the existing "no border" formatting
# (styles 4/5, e.g. row 16).
$ws.Range("A16:E16").Copy() | Out-Null
$ws.Range("A18:E18").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(18).RowHeight = 43.2

# Fill in row 18 content. Column order matches the order the new shared
# strings were appended in the source file (C, A, D, E, then the plain
# number in B).
$ws.Range("C18").Value = " Thank you, Team [team:]!"
$ws.Range("A18").Value = "SCRIPT/G01P03A/us2208.ssb"
$ws.Range("D18").Value = " Спасибо вам, Команда\n[team:]!"
$ws.Range("E18").Value = " Òðàòéáï âàí, Ëïíàîäà\n[team:]!"
$ws.Range("B18").Value = 18

# Match the view state from the diff: scrolled so row 16 is the top row,
# selection moved to C17.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("C17").Select() | Out-Null
